# Updates corresponding to the "Actualización automática 2025-06-18 14:35:09" commit.
# Adds a sale of 177.41 (PORCELANATO) for HIDALGO HIDALGO PEDRO GUSTAVO / CHONTASI SIMBAÑA SILVIA JANETH
# in June, and propagates the change through the dependent summary sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": PORCELANATO column (M), advisor row 7 ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M7").Value = 177.41
$wsVentasGrupo.Range("M22").Value = "2 de 20"

# --- Sheet "VENTA MENSUAL": junio column (F), advisor row 7 and TOTAL row 22 ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F7").Value = 177.41
$wsVentaMensual.Range("F22").Value = 168.56

# --- Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO row 16 and TOTAL row 19 ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 223.32
$wsCumplimiento.Range("E16").Value = 29309.12
$wsCumplimiento.Range("F16").Value = 0.007561854015448774

$wsCumplimiento.Range("D19").Value = 168.56
$wsCumplimiento.Range("E19").Value = 50218.63762291768
$wsCumplimiento.Range("F19").Value = 0.003345294200750184
